$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: fill in Band hours (B and D columns) ---
$ws.Range("B4").Value = 24
$ws.Range("D4").Value = 20

# --- Row 5 ---
$ws.Range("B5").Value = 14
$ws.Range("D5").Value = 12

# --- Row 6: was mostly blank, now fully populated with numbers ---
# Column B previously carried a [H]:MM:SS (time) number format; the new
# value is a plain hour count, so switch it back to General like the
# rest of the row before writing the value.
$ws.Range("B6").NumberFormat = "General"
$ws.Range("B6").Value = 26
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 20
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# --- Row 7 ---
$ws.Range("B7").NumberFormat = "General"
$ws.Range("B7").Value = 32
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 30
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0

# --- Row 8 ---
$ws.Range("B8").NumberFormat = "General"
$ws.Range("B8").Value = 24
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 22
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0

# --- Row 9 ---
$ws.Range("B9").NumberFormat = "General"
$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 10
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0

# --- New task rows 10-12: add EVF task names and estimates ---
$ws.Range("A10").Value = "EVF 1 - Map selection Multiple councils"
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 6

$ws.Range("A11").Value = "EVF 2 - Data Selection allows for Query Building"
$ws.Range("D11").Value = 20
$ws.Range("E11").Value = 10

$ws.Range("A12").Value = "EVF 3 - Enhanced Data Visualization"
$ws.Range("D12").Value = 30
$ws.Range("E12").Value = 35

# --- Selection moved to E13 ---
$ws.Range("E13").Select()
